$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 (pushes existing rows 19-50 down to 20-51),
# for the new "FAF and NASS processing" pipeline step
# (FAF/tabulate_annual_cropland.R).
$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = "FAF and NASS processing"
$ws.Cells.Item(19, 2).Value = "FAF/tabulate_annual_cropland.R"
$ws.Cells.Item(19, 3).Value = "cropland and pastureland totals joined with FAF; NASS output with imputed data"
$ws.Cells.Item(19, 4).Value = "NASS2012_receipts_workers_land_NAICS_imputed.csv; cropland_by_county_FAF_joined.csv"
$ws.Cells.Item(19, 5).Value = "cropland and pastureland totals with FAF and with annual proportion"
$ws.Cells.Item(19, 6).Value = "cropland_by_county_FAF_with_annual.csv"

$ws.Rows.Item(19).RowHeight = 45

# Update the view so the selection reflects where editing happened.
$ws.Range("G19").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 2
